$d = $word.ActiveDocument

# Locate the paragraph that starts the long bold "NLP over Machine Learning"
# commentary block - this is the paragraph whose text content must be wiped,
# while keeping the single <w:lastRenderedPageBreak/> marker that lives in
# its first run.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("The reason I picked NLP over Machine Learning")) {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -lt 0) {
    throw "Could not locate the target paragraph"
}

$target = $d.Paragraphs.Item($targetIndex)

# Insert a brand-new run containing only <w:lastRenderedPageBreak/> right at
# the very start of the paragraph. Because this element carries no text, it
# does not shift any of the existing Range offsets that follow it.
$insertionPoint = $d.Range($target.Range.Start, $target.Range.Start)
$pageBreakRunXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:lastRenderedPageBreak/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertionPoint.InsertXML($pageBreakRunXml)

# Now remove all of the original (bold) text runs that made up the rest of
# the paragraph, leaving only the paragraph mark (and the new run we just
# added, which holds no characters so it is outside this range).
$textRange = $d.Range($target.Range.Start, $target.Range.End - 1)
$textRange.Text = ""

# The paragraph that immediately follows is the lone empty paragraph the
# diff removes entirely - deleting its range merges it away, joining the
# (now pagebreak-only) paragraph directly with the "Stage 2 - ..." one.
$emptyFollower = $d.Paragraphs.Item($targetIndex + 1)
if ($emptyFollower.Range.Text.Trim().Length -eq 0) {
    $emptyFollower.Range.Delete()
}
